# Apply crypto list price/volume updates from the commit diff.
# D-column values are prefixed with a leading apostrophe (quote-prefix)
# so Excel stores them as text, matching the original inlineStr cells
# (preventing numeric auto-conversion / trailing-zero loss, e.g. 0.0890 -> 0.089).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.136.72"

$ws.Range("D3").Value = "'1.678.32"
$ws.Range("E3").Value = "  +0.01%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'214.13"
$ws.Range("E5").Value = "  -0.89%  "

$ws.Range("E6").Value = "  -0.19%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").Value = "'22.78"
$ws.Range("E8").Value = "  +6.97%  "

$ws.Range("E10").Value = "  -0.14%  "

$ws.Range("D11").Value = "'0.0890"
$ws.Range("E11").Value = "  +0.08%  "

$ws.Range("D12").Value = "'1.916.46"
$ws.Range("E12").Value = "  -0.02%  "

$ws.Range("D13").Value = "'1.676.43"
$ws.Range("E13").Value = "  -2.98%  "

$ws.Range("E14").Value = "  +2.34%  "

$ws.Range("D15").Value = "'0.551"
$ws.Range("E15").Value = "  +3.45%  "

$ws.Range("E16").Value = "  +0.20%  "

$ws.Range("D17").Value = "'27.108.23"
$ws.Range("E17").Value = "  +0.25%  "

$ws.Range("D18").Value = "'235.62"
$ws.Range("E18").Value = "  +0.12%  "

$ws.Range("D19").Value = "'7.91"
$ws.Range("E19").Value = "  -2.78%  "

$ws.Range("E20").Value = "  +0.33%  "

$ws.Range("E21").Value = "  +0.11%  "

$ws.Range("E22").Value = "  +1.91%  "

$ws.Range("D23").Value = "'9.52"
$ws.Range("E23").Value = "  +2.92%  "

$ws.Range("E24").Value = "  -1.21%  "

$ws.Range("D25").Value = "'147.23"
$ws.Range("E25").Value = "  +0.44%  "

$ws.Range("D26").Value = "'7.43"
$ws.Range("E26").Value = "  +2.61%  "

$ws.Range("D27").Value = "'16.34"
$ws.Range("E27").Value = "  -0.54%  "

$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("E30").Value = "  +0.71%  "

$ws.Range("E31").Value = "  -0.34%  "

$ws.Range("E32").Value = "  +0.11%  "

$ws.Range("D33").Value = "'1.543.19"
$ws.Range("E33").Value = "  +0.54%  "

$ws.Range("E34").Value = "  +1.72%  "

$ws.Range("E35").Value = "  -2.78%  "

$ws.Range("E36").Value = "  +2.97%  "

$ws.Range("D37").Value = "'0.939"
$ws.Range("E37").Value = "  +2.34%  "

$ws.Range("E38").Value = "  -0.33%  "

$ws.Range("E39").Value = "  -1.36%  "

$ws.Range("E40").Value = "  +2.75%  "

$ws.Range("D41").Value = "'5.79"
$ws.Range("E41").Value = "  +3.37%  "

$ws.Range("D42").Value = "'69.52"
$ws.Range("E42").Value = "  +2.35%  "

$ws.Range("E43").Value = "  +0.06%  "

$ws.Range("E44").Value = "  -0.08%  "

$ws.Range("D45").Value = "'1.823.23"

$ws.Range("D46").Value = "'0.778"
$ws.Range("E46").Value = "  -0.32%  "

$ws.Range("D47").Value = "'89.81"
$ws.Range("E47").Value = "  -0.60%  "

$ws.Range("D48").Value = "'0.0₆0111"
$ws.Range("E48").Value = "  +3.38%  "

$ws.Range("E49").Value = "  +6.03%  "

$ws.Range("D50").Value = "'8.24"
$ws.Range("E50").Value = "  +2.74%  "

$ws.Range("E51").Value = "  -0.01%  "
